$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Non-numeric-looking text columns (B, C, E) - direct value assignment is safe
$ws.Range('E2').Value = '  -0.32%  '
$ws.Range('E3').Value = '  -1.92%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('E5').Value = '  -0.39%  '
$ws.Range('E6').Value = '  -3.68%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  -8.63%  '
$ws.Range('E9').Value = '  -2.02%  '
$ws.Range('E10').Value = '  -1.56%  '
$ws.Range('E11').Value = '  +1.49%  '
$ws.Range('E12').Value = '  -1.47%  '
$ws.Range('E13').Value = '  -1.12%  '
$ws.Range('E14').Value = '  -3.92%  '
$ws.Range('E15').Value = '  -1.93%  '
$ws.Range('E16').Value = '  +0.27%  '
$ws.Range('E17').Value = '  -3.99%  '
$ws.Range('E18').Value = '  -1.94%  '
$ws.Range('E19').Value = '  -0.44%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('E20').Value = '  -0.16%  '
$ws.Range('B21').Value = 'Polkadot'
$ws.Range('C21').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('E21').Value = '  -0.29%  '
$ws.Range('E22').Value = '  -3.41%  '
$ws.Range('E23').Value = '  +0.06%  '
$ws.Range('E24').Value = '  -0.16%  '
$ws.Range('E25').Value = '  -7.93%  '
$ws.Range('E26').Value = '  +7.13%  '
$ws.Range('E27').Value = '  +0.13%  '
$ws.Range('E28').Value = '  -2.08%  '
$ws.Range('E29').Value = '  -1.71%  '
$ws.Range('E30').Value = '  -5.43%  '
$ws.Range('E31').Value = '  -5.18%  '
$ws.Range('B32').Value = 'Kaspa'
$ws.Range('C32').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('E32').Value = '  +0.14%  '
$ws.Range('B33').Value = 'PEPE'
$ws.Range('C33').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('E33').Value = '  -9.52%  '
$ws.Range('E34').Value = '  -2.69%  '
$ws.Range('E35').Value = '  -5.12%  '
$ws.Range('E36').Value = '  -0.03%  '
$ws.Range('B37').Value = 'PolygonEcosystemToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('E37').Value = '  -0.68%  '
$ws.Range('B38').Value = 'NEARProtocol'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('E38').Value = '  -2.50%  '
$ws.Range('E39').Value = '  +1.51%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('E40').Value = '  -5.22%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('E41').Value = '  +4.04%  '
$ws.Range('E42').Value = '  +3.80%  '
$ws.Range('E44').Value = '  +0.76%  '
$ws.Range('E45').Value = '  +1.05%  '
$ws.Range('E46').Value = '  -0.68%  '
$ws.Range('E47').Value = '  -10.85%  '
$ws.Range('E48').Value = '  -0.08%  '
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('E49').Value = '  -2.56%  '
$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('E50').Value = '  -7.77%  '
$ws.Range('E51').Value = '  -3.06%  '

# Column D (price) values look numeric; force text storage, then restore default style
$priceCells = [ordered]@{
    'D2' = '60.671.06'
    'D3' = '2.335.96'
    'D5' = '541.83'
    'D6' = '135.55'
    'D8' = '0.523'
    'D9' = '2.336.51'
    'D12' = '5.27'
    'D13' = '0.340'
    'D14' = '24.47'
    'D15' = '2.760.22'
    'D16' = '60.974.73'
    'D17' = '0.0000158'
    'D18' = '2.337.79'
    'D19' = '10.56'
    'D20' = '317.30'
    'D21' = '4.09'
    'D22' = '6.49'
    'D24' = '63.10'
    'D25' = '1.66'
    'D26' = '8.36'
    'D28' = '2.451.18'
    'D29' = '7.88'
    'D30' = '491.71'
    'D31' = '1.36'
    'D32' = '0.145'
    'D33' = '0.0₃0851'
    'D35' = '1.48'
    'D36' = '1.00'
    'D37' = '0.374'
    'D38' = '4.54'
    'D39' = '18.39'
    'D40' = '5.20'
    'D41' = '1.80'
    'D42' = '142.61'
    'D44' = '40.53'
    'D45' = '141.64'
    'D46' = '3.52'
    'D47' = '2.01'
    'D48' = '0.0514'
    'D49' = '0.565'
    'D50' = '18.85'
    'D51' = '0.0895'
}
foreach ($addr in $priceCells.Keys) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $priceCells[$addr]
    $rng.Style = "Normal"
}